$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove all existing comments first. This engine does NOT shift comment
#    anchors when rows are inserted (unlike real Excel), so the safest way to
#    end up with comments on the correct (post-insert) cells is to drop them
#    all up front, do the row surgery, then re-create every comment (both the
#    ones that already existed and the two brand new ones) at their final
#    co-ordinates.
# ---------------------------------------------------------------------------
$savedComments = New-Object System.Collections.ArrayList
foreach ($cm in $ws.Comments) {
    $entry = @{ Ref = $cm.Parent.Address(); Text = $cm.Text() }
    [void]$savedComments.Add($entry)
}
while ($ws.Comments.Count -gt 0) {
    $ws.Comments.Item(1).Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert the two new rows.
#    - row 14 : Polymaker PolyMax PC   (pushes the old 14..30 down to 15..31)
#    - row 31 (after the first insert, the old last row 30 now sits at 31) :
#      Taulman PCTPE                   (pushes that row down to 32)
# ---------------------------------------------------------------------------
$ws.Rows("14:14").Insert()
$ws.Rows("31:31").Insert()

# New row 14: Polymaker PolyMax PC
$ws.Range("A14").Value2 = "Polymaker PolyMax PC"
$ws.Range("B14").Value2 = "Taulman Alloy 910"
$ws.Range("C14").Value2 = "Polymaker PC-Max fans"
$ws.Range("D14").Value2 = "Original Prusa i3 MK3 purgebubble"
$ws.Range("E14").Value2 = "Yes"

# New row 31: Taulman PCTPE
$ws.Range("A31").Value2 = "Taulman PCTPE"
$ws.Range("B31").Value2 = "Taulman Alloy 910"
$ws.Range("C31").Value2 = "Taulman PCTPE"
$ws.Range("D31").Value2 = "Original Prusa i3 MK3 purgebubble"
$ws.Range("E31").Value2 = "Yes"

# ---------------------------------------------------------------------------
# 3) Bring formatting (borders/fills) on the two brand-new rows into line
#    with their neighbours, same as Excel does when you insert a copied row.
# ---------------------------------------------------------------------------
$ws.Range("A15:F15").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)

$ws.Range("A30:F30").Copy()
$ws.Range("A31:F31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Re-create every comment. Existing comments whose row was at/after the
#    first insertion point (14) shift down by one; those at/after the second
#    insertion point (31, measured after the first shift) shift down by one
#    more. Plus the two new comments on D14 and D31.
# ---------------------------------------------------------------------------
function Shift-Row([int]$row) {
    if ($row -ge 14) { $row = $row + 1 }
    if ($row -ge 31) { $row = $row + 1 }
    return $row
}

foreach ($entry in $savedComments) {
    $ref = $entry.Ref -replace '\$', ''
    if ($ref -match '^([A-Z]+)([0-9]+)$') {
        $col = $matches[1]
        $row = [int]$matches[2]
        $newRow = Shift-Row $row
        $newRef = "$col$newRow"
        $ws.Range($newRef).AddComment($entry.Text) | Out-Null
    }
}

$ws.Range("D14").AddComment("Darragh Broadbent:`nPrint quality improved by very conservative use of the cooling fan, especially small detail and short layer time, printed with a brim, minimal if any warping.") | Out-Null
$ws.Range("D31").AddComment("Darragh Broadbent:`nNot fantastic quality, stringing and warping artifacts, printed with a brim, heavy warping.") | Out-Null

# ---------------------------------------------------------------------------
# 5) Restore the selection shown in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("J27").Select()
